$d = $word.ActiveDocument

# --- Change 1: merge the "{{#each exam.users}}" paragraph's runs into one run ---
$d.Content.Find.Execute("{{#each exam.users}}", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "{{#each exam.users}}", 2) | Out-Null

# --- Change 2: merge the "{{@index}}. User: {{name}} has:" paragraph's runs into one run ---
$d.Content.Find.Execute("{{@index}}. User: {{name}} has:", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "{{@index}}. User: {{name}} has:", 2) | Out-Null

# --- Change 3: fix the paragraph-mark formatting of the final "{{/each}}" paragraph, then
#     insert two new paragraphs (exam.options each-loop) right before it ---

# Locate the paragraph that still carries the old bold/32pt paragraph-mark formatting
# and contains the text "{{/each}}" (the one that closes the exam.users loop).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "{{/each}}") {
        $target = $p
    }
}

# Normalize its paragraph-mark run properties: drop bold, shrink to 28/28 (matches the
# run formatting already used by the text inside it).
$xml = $target.Range.WordOpenXML
$xmlFixed = $xml.Replace('<w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/>', '<w:sz w:val="28"/><w:szCs w:val="28"/>')
$target.Range.InsertXML($xmlFixed) | Out-Null

# Re-resolve the (now reformatted) paragraph and its index so we can insert after it.
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.TrimEnd([char]13, [char]7) -eq "{{/each}}" -and $cand.Range.Start -eq $target.Range.Start) {
        $targetIndex = $i
    }
}

# Insert two brand-new empty paragraphs right after it; they inherit the (now-fixed)
# paragraph-mark formatting automatically.
$d.Paragraphs.Item($targetIndex).Range.InsertParagraphAfter() | Out-Null
$d.Paragraphs.Item($targetIndex + 1).Range.InsertParagraphAfter() | Out-Null

# First new paragraph: "{{#each exam.options}}" (single run).
$d.Paragraphs.Item($targetIndex + 1).Range.Text = "{{#each exam.options}}"

# Second new paragraph: two runs, "{{@index}}" and ". {{value}}".
$snippet = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
    '<w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>{{@index}}</w:t></w:r>' + `
    '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/><w:lang w:val="en-US"/></w:rPr><w:t>. {{value}}</w:t></w:r>' + `
    '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs.Item($targetIndex + 2).Range.InsertXML($snippet) | Out-Null
